$wb = $excel.ActiveWorkbook

# Update the "Ready for handoff" status cells to "In Translation" across
# the Overview sheet (columns E:F) and the per-language sheets (column C).
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value = "In Translation"
